$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.491.16"
Set-TextValue $ws.Range("E2") "  +0.18%  "
Set-TextValue $ws.Range("D3") "1.858.11"
Set-TextValue $ws.Range("E3") "  +0.51%  "
Set-TextValue $ws.Range("D4") "1.004"
Set-TextValue $ws.Range("E4") "  +0.29%  "
Set-TextValue $ws.Range("E5") "  +0.43%  "
Set-TextValue $ws.Range("D6") "0.6333"
Set-TextValue $ws.Range("E6") "  +1.12%  "
Set-TextValue $ws.Range("E7") "  +0.17%  "
Set-TextValue $ws.Range("D8") "0.07583"
Set-TextValue $ws.Range("E8") "  -1.10%  "
Set-TextValue $ws.Range("D9") "0.2927"
Set-TextValue $ws.Range("E9") "  +0.41%  "
Set-TextValue $ws.Range("D10") "24.60"
Set-TextValue $ws.Range("E10") "  -0.60%  "
Set-TextValue $ws.Range("D11") "0.07767"
Set-TextValue $ws.Range("E11") "  +0.21%  "
Set-TextValue $ws.Range("D12") "1.857.12"
Set-TextValue $ws.Range("E12") "  -0.60%  "
Set-TextValue $ws.Range("D13") "5.043"
Set-TextValue $ws.Range("E13") "  +0.32%  "
Set-TextValue $ws.Range("D14") "0.6853"
Set-TextValue $ws.Range("E14") "  +0.63%  "
Set-TextValue $ws.Range("D15") "0.00001051"
Set-TextValue $ws.Range("E15") "  -1.82%  "
Set-TextValue $ws.Range("D16") "83.49"
Set-TextValue $ws.Range("E16") "  +0.00%  "
Set-TextValue $ws.Range("D17") "2.118.66"
Set-TextValue $ws.Range("E17") "  +0.62%  "
Set-TextValue $ws.Range("D18") "6.156"
Set-TextValue $ws.Range("E18") "  -0.16%  "
Set-TextValue $ws.Range("D19") "29.498.92"
Set-TextValue $ws.Range("E19") "  +0.11%  "
Set-TextValue $ws.Range("D20") "230.36"
Set-TextValue $ws.Range("E20") "  +0.92%  "
Set-TextValue $ws.Range("D21") "12.40"
Set-TextValue $ws.Range("E21") "  +0.16%  "
Set-TextValue $ws.Range("E22") "  +0.21%  "
Set-TextValue $ws.Range("D23") "7.519"
Set-TextValue $ws.Range("E23") "  +1.39%  "
Set-TextValue $ws.Range("D24") "1.004"
Set-TextValue $ws.Range("E24") "  +0.20%  "
Set-TextValue $ws.Range("D25") "159.56"
Set-TextValue $ws.Range("E25") "  +1.21%  "
Set-TextValue $ws.Range("D26") "0.1398"
Set-TextValue $ws.Range("E26") "  +1.93%  "
Set-TextValue $ws.Range("D27") "8.470"
Set-TextValue $ws.Range("E27") "  +0.89%  "
Set-TextValue $ws.Range("D28") "17.74"
Set-TextValue $ws.Range("E28") "  +0.18%  "
Set-TextValue $ws.Range("D29") "1.422"
Set-TextValue $ws.Range("E29") "  +5.29%  "
Set-TextValue $ws.Range("D30") "1.483"
Set-TextValue $ws.Range("E30") "  +1.17%  "
Set-TextValue $ws.Range("D31") "0.05708"
Set-TextValue $ws.Range("E31") "  +1.43%  "
Set-TextValue $ws.Range("D32") "4.157"
Set-TextValue $ws.Range("E32") "  +1.00%  "
Set-TextValue $ws.Range("D33") "4.079"
Set-TextValue $ws.Range("E33") "  +1.19%  "
Set-TextValue $ws.Range("D34") "1.834"
Set-TextValue $ws.Range("E34") "  -0.40%  "
Set-TextValue $ws.Range("D35") "1.159"
Set-TextValue $ws.Range("E35") "  -0.13%  "
Set-TextValue $ws.Range("D36") "0.6969"
Set-TextValue $ws.Range("E36") "  -0.90%  "
Set-TextValue $ws.Range("E37") "  +0.03%  "
Set-TextValue $ws.Range("D38") "1.258.14"
Set-TextValue $ws.Range("E38") "  +2.36%  "
Set-TextValue $ws.Range("D39") "0.01836"
Set-TextValue $ws.Range("E39") "  +2.67%  "
Set-TextValue $ws.Range("D40") "2.783"
Set-TextValue $ws.Range("E40") "  +0.58%  "
Set-TextValue $ws.Range("D41") "6.518"
Set-TextValue $ws.Range("E41") "  -0.30%  "
Set-TextValue $ws.Range("D42") "0.9086"
Set-TextValue $ws.Range("E42") "  +0.52%  "
Set-TextValue $ws.Range("E43") "  +0.09%  "
Set-TextValue $ws.Range("D44") "2.021.06"
Set-TextValue $ws.Range("E44") "  +0.39%  "
Set-TextValue $ws.Range("D45") "101.71"
Set-TextValue $ws.Range("E45") "  -0.07%  "
Set-TextValue $ws.Range("D46") "66.24"
Set-TextValue $ws.Range("E46") "  +0.57%  "
Set-TextValue $ws.Range("D47") "7.155"
Set-TextValue $ws.Range("D48") "0.1170"
Set-TextValue $ws.Range("E48") "  +1.54%  "
Set-TextValue $ws.Range("D49") "9.064"
Set-TextValue $ws.Range("E49") "  +0.69%  "
Set-TextValue $ws.Range("D50") "0.3981"
Set-TextValue $ws.Range("E51") "  +0.51%  "
